$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionsText = @'
questions = [
    {
        "title": "You need to send an email to which Alex needs to respond, while Kelly only needs to be informed of the conversation.How should you address the email?",
        "ques_type": 2,
        "options": [
            "Put Alex in the Bcc field and Kelly in the Cc field.",
            "Put Alex in the To field and Kelly in the Cc field.",
            "Put Alex in the Cc field and Kelly in the Bcc field.",
            "Put Alex in the Cc field and Kelly in the To field."
        ],
        "score": "Put Alex in the To field and Kelly in the Cc field."
    },
    {
        "title": "You received an email where you are Cc'd. The email contains a complex issue that needs to be addressed, but you're not directly asked for a solution.What should be your response etiquette in this situation?",
        "ques_type": 2,
        "options": [
            "Forward the email to another team member for input before replying. ",
            "Reply, asking for clarification on your role regarding the issue.",
            "Reply immediately, promising to look into the issue.",
            "Reply if you can provide insightful input or a potential solution."
        ],
        "score": "Reply if you can provide insightful input or a potential solution."
    },
    {
        "title": "You emailed some coworkers with whom you are working on a project. After sending the email, you realize that you forgot to include two of the team members on the recipient list. What should you do?",
        "ques_type": 2,
        "options": [
            "Select the message in the Sent Items folder, navigate to the Message tab, then select Actions &gt Resend This Message.",
            "Select the message in the Sent Items folder, navigate to the Message tab, then select Actions &gt Recall This Message.",
            "Navigate to the Outbox folder, press Recover items deleted from this folder, find the sent message, and select Restore.",
            "Select the message in the Sent Items folder, navigate to the Message tab, then select Actions &gt Edit Message."
        ],
        "score": "Select the message in the Sent Items folder, navigate to the Message tab, then select Actions &gt Resend This Message."
    },
    {
        "title": "You set up a departmental meeting using Microsoft Outlook with default settings. However, you realize you invited someone by mistake and want to delete them as an attendee. You open the original meeting invitation and delete the person's name from the invite list.What should you do to ensure the correct people are notified?",
        "ques_type": 2,
        "options": [
            "Look for a dialog box that will guide your next step.",
            "Send an update to all attendees immediately.",
            "Inform your supervisor of the change. ",
            "Cancel the meeting and recreate it from scratch. "
        ],
        "score": "Look for a dialog box that will guide your next step."
    }
]
'@

$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $questionsText
$ws.Rows.Item(1).AutoFit()
$ws.Range("A2").ClearContents()

Write-Output "done"
